# Add the "date stamp" cell (About!C1) that accompanies the title in A1.
# Source XML stores this as the raw serial number 44307 (2021-04-21) with
# a date-formatted style (numFmtId 14, the builtin short-date format), so
# assign a plain numeric value and a matching NumberFormat rather than a
# .NET DateTime (which would otherwise carry a fractional time-of-day part).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")

$ws.Range("C1").Value = 44307
$ws.Range("C1").NumberFormat = "mm-dd-yy"
